$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 2499
$ws.Cells.Item(51, 9).Value = 2000
$ws.Cells.Item(51, 10).Value = 2998
$ws.Cells.Item(51, 11).Value = 2000
$ws.Cells.Item(51, 12).Value = 2998
$ws.Cells.Item(51, 13).Value = -1516
$ws.Cells.Item(51, 14).Value = -3966

$ws.Cells.Item(94, 8).Value = 723.5
$ws.Cells.Item(94, 9).Value = 631.3333
$ws.Cells.Item(94, 10).Value = 1000
$ws.Cells.Item(94, 11).Value = 631.3333
$ws.Cells.Item(94, 12).Value = 1000
$ws.Cells.Item(94, 13).Value = -180.3333
$ws.Cells.Item(94, 14).Value = -1902

$ws.Cells.Item(112, 8).Value = 1788.375
$ws.Cells.Item(112, 9).Value = 1152.3334
$ws.Cells.Item(112, 10).Value = 2170
$ws.Cells.Item(112, 11).Value = 3457.0002
$ws.Cells.Item(112, 12).Value = 6510
$ws.Cells.Item(112, 13).Value = -2349.0002
$ws.Cells.Item(112, 14).Value = -8726

$ws.Cells.Item(125, 8).Value = 300001100
$ws.Cells.Item(125, 10).Value = 125001380
$ws.Cells.Item(125, 12).Value = 1125012420
$ws.Cells.Item(125, 14).Value = -1125017340

$ws.Cells.Item(137, 8).Value = 1216.3334
$ws.Cells.Item(137, 9).Value = 931
$ws.Cells.Item(137, 11).Value = 2793
$ws.Cells.Item(137, 13).Value = -243

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 984.3333
$ws.Cells.Item(102, 9).Value = 982.2857
$ws.Cells.Item(102, 10).Value = 991.5
$ws.Cells.Item(102, 11).Value = 982.2857
$ws.Cells.Item(102, 12).Value = 991.5
$ws.Cells.Item(102, 13).Value = 639.7143
$ws.Cells.Item(102, 14).Value = -4235.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2124.5
$ws.Cells.Item(86, 9).Value = 1987
$ws.Cells.Item(86, 10).Value = 2399.5
$ws.Cells.Item(86, 11).Value = 1987
$ws.Cells.Item(86, 12).Value = 2399.5
$ws.Cells.Item(86, 13).Value = -864
$ws.Cells.Item(86, 14).Value = -4645.5

$ws.Cells.Item(89, 8).Value = 2124.5
$ws.Cells.Item(89, 9).Value = 1987
$ws.Cells.Item(89, 10).Value = 2399.5
$ws.Cells.Item(89, 11).Value = 9935
$ws.Cells.Item(89, 12).Value = 11997.5
$ws.Cells.Item(89, 13).Value = -4319
$ws.Cells.Item(89, 14).Value = -23229.5

$ws.Cells.Item(105, 8).Value = 2870.9375
$ws.Cells.Item(105, 10).Value = 3792.75
$ws.Cells.Item(105, 12).Value = 3792.75
$ws.Cells.Item(105, 14).Value = -7286.75

$ws.Cells.Item(107, 8).Value = 3281.6
$ws.Cells.Item(107, 9).Value = 2987.875
$ws.Cells.Item(107, 11).Value = 2987.875
$ws.Cells.Item(107, 13).Value = -1067.875

$ws.Cells.Item(134, 8).Value = 436.33334
$ws.Cells.Item(134, 9).Value = 436.33334
$ws.Cells.Item(134, 11).Value = 1309.00002
$ws.Cells.Item(134, 13).Value = 1225.99998

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(19, 8).Value = 210
$ws.Cells.Item(19, 9).Value = 210
$ws.Cells.Item(19, 11).Value = 210
$ws.Cells.Item(19, 13).Value = -40

$ws.Cells.Item(22, 8).Value = 29514.666
$ws.Cells.Item(22, 9).Value = 1725.6
$ws.Cells.Item(22, 10).Value = 64251
$ws.Cells.Item(22, 11).Value = 1725.6
$ws.Cells.Item(22, 12).Value = 64251
$ws.Cells.Item(22, 13).Value = -1375.6
$ws.Cells.Item(22, 14).Value = -64951

$ws.Cells.Item(24, 8).Value = 210
$ws.Cells.Item(24, 9).Value = 210
$ws.Cells.Item(24, 11).Value = 210
$ws.Cells.Item(24, 13).Value = -40

$ws.Cells.Item(31, 8).Value = 2883.25
$ws.Cells.Item(31, 9).Value = 2363.1667
$ws.Cells.Item(31, 11).Value = 2363.1667
$ws.Cells.Item(31, 13).Value = -2068.1667

$ws.Cells.Item(34, 8).Value = 2883.25
$ws.Cells.Item(34, 9).Value = 2363.1667
$ws.Cells.Item(34, 11).Value = 2363.1667
$ws.Cells.Item(34, 13).Value = -2161.1667

$ws.Cells.Item(100, 8).Value = 150384.5
$ws.Cells.Item(100, 10).Value = 150384.5
$ws.Cells.Item(100, 12).Value = 150384.5
$ws.Cells.Item(100, 14).Value = -152548.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 55000000
$ws.Cells.Item(4, 9).Value = 55000000
$ws.Cells.Item(4, 11).Value = 165000000
$ws.Cells.Item(4, 13).Value = -164999888

$ws.Cells.Item(60, 8).Value = 845.75
$ws.Cells.Item(60, 9).Value = 496.33334
$ws.Cells.Item(60, 11).Value = 1489.00002
$ws.Cells.Item(60, 13).Value = -1238.00002

$ws.Cells.Item(94, 8).Value = 14240.5
$ws.Cells.Item(94, 9).Value = 4808
$ws.Cells.Item(94, 11).Value = 14424
$ws.Cells.Item(94, 13).Value = -13748

$ws.Cells.Item(131, 8).Value = 2420.1177
$ws.Cells.Item(131, 9).Value = 1488.625
$ws.Cells.Item(131, 10).Value = 3248.111
$ws.Cells.Item(131, 11).Value = 4465.875
$ws.Cells.Item(131, 12).Value = 9744.332999999999
$ws.Cells.Item(131, 13).Value = 574.125
$ws.Cells.Item(131, 14).Value = -19824.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(36, 8).Value = 10459
$ws.Cells.Item(36, 9).Value = 10459
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 10459
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = -9974
$ws.Cells.Item(36, 14).ClearContents()

$ws.Cells.Item(113, 8).Value = 1534.5
$ws.Cells.Item(113, 9).Value = 1534.5
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 1534.5
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = 635.5
$ws.Cells.Item(113, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 2381.25
$ws.Cells.Item(122, 9).Value = 1532.5
$ws.Cells.Item(122, 11).Value = 4597.5
$ws.Cells.Item(122, 13).Value = -2147.5

$ws.Cells.Item(126, 8).Value = 3548.5334
$ws.Cells.Item(126, 10).Value = 3874.625
$ws.Cells.Item(126, 12).Value = 11623.875
$ws.Cells.Item(126, 14).Value = -16563.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 8194.65
$ws.Cells.Item(7, 9).Value = 4303.5
$ws.Cells.Item(7, 10).Value = 8627
$ws.Cells.Item(7, 11).Value = 4303.5
$ws.Cells.Item(7, 12).Value = 8627
$ws.Cells.Item(7, 13).Value = -4191.5
$ws.Cells.Item(7, 14).Value = -8851

$ws.Cells.Item(93, 8).Value = 9000
$ws.Cells.Item(93, 9).Value = 9000
$ws.Cells.Item(93, 11).Value = 9000
$ws.Cells.Item(93, 13).Value = -7752

$ws.Cells.Item(103, 8).Value = 18500
$ws.Cells.Item(103, 10).Value = 18500
$ws.Cells.Item(103, 12).Value = 18500
$ws.Cells.Item(103, 14).Value = -20844

$ws.Cells.Item(126, 8).Value = 8194.65
$ws.Cells.Item(126, 9).Value = 4303.5
$ws.Cells.Item(126, 10).Value = 8627
$ws.Cells.Item(126, 11).Value = 12910.5
$ws.Cells.Item(126, 12).Value = 25881
$ws.Cells.Item(126, 13).Value = -10440.5
$ws.Cells.Item(126, 14).Value = -30821

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(29, 8).Value = 2400
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 2400
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 2400
$ws.Cells.Item(29, 14).Value = -2980
$ws.Cells.Item(29, 13).ClearContents()

$ws.Cells.Item(43, 8).Value = 20000
$ws.Cells.Item(43, 9).Value = 20000
$ws.Cells.Item(43, 11).Value = 20000
$ws.Cells.Item(43, 13).Value = -19851

$ws.Cells.Item(113, 8).Value = 498.2
$ws.Cells.Item(113, 10).Value = 999.75
$ws.Cells.Item(113, 12).Value = 2999.25
$ws.Cells.Item(113, 14).Value = -7339.25

$ws.Cells.Item(126, 8).Value = 4222.7646
$ws.Cells.Item(126, 10).Value = 5800
$ws.Cells.Item(126, 12).Value = 17400
$ws.Cells.Item(126, 14).Value = -22340
